$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<sere>"
$ws.Range("C2").Value = 39

$ws.Range("C3").Value = 34

$ws.Range("B4").Value = "<thetem>"
$ws.Range("C4").Value = 35

$ws.Range("C5").Value = 38

$ws.Range("B6").Value = "<pare>"
$ws.Range("C6").Value = 41

$ws.Range("C7").Value = 42

$ws.Range("C8").Value = 41

$ws.Range("C9").Value = 32

$ws.Range("B10").Value = "<soe>"
$ws.Range("C10").Value = 42

$ws.Range("B12").Value = "<be>"

$ws.Range("B13").Value = "<hid>"
$ws.Range("C13").Value = 40

$ws.Range("C14").Value = 36

$ws.Range("B15").Value = "<of>"
$ws.Range("C15").Value = 35
